$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2471.4285
$ws.Range("I28").Value = 2471.4285
$ws.Range("K28").Value = 2471.4285
$ws.Range("M28").Value = -1986.4285
$ws.Range("H58").Value = 799.8889
$ws.Range("H62").Value = 4166.4546
$ws.Range("I62").Value = 4229.5
$ws.Range("J62").Value = 3998.3333
$ws.Range("K62").Value = 4229.5
$ws.Range("L62").Value = 3998.3333
$ws.Range("M62").Value = -3605.5
$ws.Range("N62").Value = -5246.3333
$ws.Range("H65").Value = 4166.4546
$ws.Range("I65").Value = 4229.5
$ws.Range("J65").Value = 3998.3333
$ws.Range("K65").Value = 21147.5
$ws.Range("L65").Value = 19991.6665
$ws.Range("M65").Value = -18027.5
$ws.Range("N65").Value = -26231.6665
$ws.Range("H135").Value = 17863964
$ws.Range("I135").Value = 26317370
$ws.Range("K135").Value = 236856330
$ws.Range("M135").Value = -236853795
$ws.Range("H137").Value = 3872.25
$ws.Range("I137").Value = 3307.5
$ws.Range("K137").Value = 9922.5
$ws.Range("M137").Value = -7372.5
$ws.Range("H138").Value = 2396.6743
$ws.Range("I138").Value = 1404.1666
$ws.Range("J138").Value = 3650.3684
$ws.Range("K138").Value = 4212.4998
$ws.Range("L138").Value = 10951.1052
$ws.Range("M138").Value = 927.5002000000004
$ws.Range("N138").Value = -21231.1052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13379.615
$ws.Range("I2").Value = 20728.562
$ws.Range("K2").Value = 20728.562
$ws.Range("M2").Value = -20615.562
$ws.Range("H61").Value = 3646.7
$ws.Range("I61").Value = 1340.4166
$ws.Range("K61").Value = 1340.4166
$ws.Range("M61").Value = -1128.4166
$ws.Range("H74").Value = 2653.375
$ws.Range("I74").Value = 2630.2666
$ws.Range("K74").Value = 2630.2666
$ws.Range("M74").Value = -1756.2666
$ws.Range("H77").Value = 2653.375
$ws.Range("I77").Value = 2630.2666
$ws.Range("K77").Value = 13151.333
$ws.Range("M77").Value = -8783.332999999999
$ws.Range("H116").Value = 13379.615
$ws.Range("I116").Value = 20728.562
$ws.Range("K116").Value = 20728.562
$ws.Range("M116").Value = -18434.562
$ws.Range("H132").Value = 3449.6667
$ws.Range("I132").Value = 2962.125
$ws.Range("K132").Value = 8886.375
$ws.Range("M132").Value = -6356.375
$ws.Range("H136").Value = 3646.7
$ws.Range("I136").Value = 1340.4166
$ws.Range("K136").Value = 4021.2498
$ws.Range("M136").Value = -1471.2498
$ws.Range("H139").Value = 80130
$ws.Range("J139").Value = 80130
$ws.Range("L139").Value = 80130
$ws.Range("N139").Value = -90410

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13379.615
$ws.Range("I3").Value = 20728.562
$ws.Range("K3").Value = 20728.562
$ws.Range("M3").Value = -20614.562
$ws.Range("H70").Value = 49998.57
$ws.Range("J70").Value = 49998.57
$ws.Range("L70").Value = 49998.57
$ws.Range("N70").Value = -50584.57
$ws.Range("H73").Value = 49998.57
$ws.Range("J73").Value = 49998.57
$ws.Range("L73").Value = 49998.57
$ws.Range("N73").Value = -52026.57

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 960.9091
$ws.Range("I16").Value = 774.875
$ws.Range("J16").Value = 1457
$ws.Range("K16").Value = 774.875
$ws.Range("L16").Value = 1457
$ws.Range("M16").Value = -487.875
$ws.Range("N16").Value = -2031
$ws.Range("H31").Value = 2510.5217
$ws.Range("I31").Value = 1166.2
$ws.Range("K31").Value = 1166.2
$ws.Range("M31").Value = -871.2
$ws.Range("H34").Value = 2510.5217
$ws.Range("I34").Value = 1166.2
$ws.Range("K34").Value = 1166.2
$ws.Range("M34").Value = -964.2
$ws.Range("H58").Value = 1979.3
$ws.Range("I58").Value = 1310.3334
$ws.Range("J58").Value = 8000
$ws.Range("K58").Value = 1310.3334
$ws.Range("L58").Value = 8000
$ws.Range("M58").Value = -1107.3334
$ws.Range("N58").Value = -8406
$ws.Range("H99").Value = 5105127
$ws.Range("I99").Value = 2038268.6
$ws.Range("J99").Value = 6255199
$ws.Range("K99").Value = 2038268.6
$ws.Range("L99").Value = 6255199
$ws.Range("M99").Value = -2036770.6
$ws.Range("N99").Value = -6258195
$ws.Range("H113").Value = 960.9091
$ws.Range("I113").Value = 774.875
$ws.Range("J113").Value = 1457
$ws.Range("K113").Value = 774.875
$ws.Range("L113").Value = 1457
$ws.Range("M113").Value = 1395.125
$ws.Range("N113").Value = -5797
$ws.Range("H122").Value = 409614.47
$ws.Range("I122").Value = 426390.12
$ws.Range("K122").Value = 1279170.36
$ws.Range("M122").Value = -1276720.36
$ws.Range("H126").Value = 5105127
$ws.Range("I126").Value = 2038268.6
$ws.Range("J126").Value = 6255199
$ws.Range("K126").Value = 6114805.800000001
$ws.Range("L126").Value = 18765597
$ws.Range("M126").Value = -6112335.800000001
$ws.Range("N126").Value = -18770537
$ws.Range("H134").Value = 4159.0454
$ws.Range("J134").Value = 7095.75
$ws.Range("L134").Value = 21287.25
$ws.Range("N134").Value = -26357.25
$ws.Range("H136").Value = 1979.3
$ws.Range("I136").Value = 1310.3334
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 3931.0002
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -1381.0002
$ws.Range("N136").Value = -29100

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 262.26315
$ws.Range("I107").Value = 295.2857
$ws.Range("K107").Value = 295.2857
$ws.Range("M107").Value = 1624.7143
$ws.Range("H113").Value = 5085.6875
$ws.Range("I113").Value = 2537.2
$ws.Range("J113").Value = 9333.166999999999
$ws.Range("K113").Value = 2537.2
$ws.Range("L113").Value = 9333.166999999999
$ws.Range("M113").Value = -367.1999999999998
$ws.Range("N113").Value = -13673.167
$ws.Range("H126").Value = 5026.2
$ws.Range("I126").Value = 2313.1667
$ws.Range("K126").Value = 6939.500100000001
$ws.Range("M126").Value = -4469.500100000001
$ws.Range("H132").Value = 4930.28
$ws.Range("J132").Value = 4617.2
$ws.Range("L132").Value = 13851.6
$ws.Range("N132").Value = -18911.6
$ws.Range("H135").Value = 60000
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("N135").Value = -70140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5670.4
$ws.Range("J7").Value = 5670.4
$ws.Range("L7").Value = 5670.4
$ws.Range("N7").Value = -5894.4
$ws.Range("H16").Value = 102.77778
$ws.Range("I16").Value = 102.77778
$ws.Range("K16").Value = 102.77778
$ws.Range("M16").Value = 67.22221999999999
$ws.Range("H40").Value = 8960.593000000001
$ws.Range("I40").Value = 10985.583
$ws.Range("K40").Value = 10985.583
$ws.Range("M40").Value = -10849.583
$ws.Range("H68").Value = 6581.125
$ws.Range("I68").Value = 2759.8
$ws.Range("J68").Value = 8318.091
$ws.Range("K68").Value = 2759.8
$ws.Range("L68").Value = 8318.091
$ws.Range("M68").Value = -2010.8
$ws.Range("N68").Value = -9816.091
$ws.Range("H71").Value = 6581.125
$ws.Range("I71").Value = 2759.8
$ws.Range("J71").Value = 8318.091
$ws.Range("K71").Value = 13799
$ws.Range("L71").Value = 41590.455
$ws.Range("M71").Value = -10055
$ws.Range("N71").Value = -49078.455
$ws.Range("H126").Value = 5670.4
$ws.Range("J126").Value = 5670.4
$ws.Range("L126").Value = 17011.2
$ws.Range("N126").Value = -21951.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("H122").Value = 4403.0293
$ws.Range("I122").Value = 2122.2222
$ws.Range("J122").Value = 13200.429
$ws.Range("K122").Value = 6366.6666
$ws.Range("L122").Value = 39601.287
$ws.Range("M122").Value = -3916.6666
$ws.Range("N122").Value = -44501.287
$ws.Range("H126").Value = 2491.077
$ws.Range("I126").Value = 2198.4
$ws.Range("K126").Value = 6595.200000000001
$ws.Range("M126").Value = -4125.200000000001
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("M129").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("M140").ClearContents()
